$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.817.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.501.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.503.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.944.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.648.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.498.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.409"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0755"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.831"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "273.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.34%  "
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.592"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0931"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.753.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.30%  "
